# Applies the cryptos.xlsx price/volume update described in the commit
# "Updated cryptos list on Wed May 31 03:49:23 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. the thousand-separated
# '27.676.38') that must stay literal text. Force Text format before writing
# so Excel doesn't auto-coerce them to numbers, then restore the Normal style
# so the cell's formatting matches the original (unstyled) cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.676.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.896.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5253"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.69%  "
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07241"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9022"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07634"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.878.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.434"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9994"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008664"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.697.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.144"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.136.93"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.597"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.857"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.177"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.837"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.820"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09152"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.157"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("E35").Value = "  -1.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7733"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.570"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("E39").Value = "  +2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5563"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.090"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.700"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "117.47"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.719"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1510"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4803"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9991"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.593"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "37.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
